$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: add Time (C14), extend Goal (D14) with a 3rd line, add Comment (E14)
$ws.Cells.Item(14, 3).Value = "08:22`n18:01"
$ws.Cells.Item(14, 3).HorizontalAlignment = -4108
$ws.Cells.Item(14, 3).VerticalAlignment = -4108
$ws.Cells.Item(14, 3).WrapText = $true

$ws.Cells.Item(14, 4).Value = "1. 앱 화면별 Sequence 작성`n2. 앱 화면별 어떻게 구성할 것인지 구상`n3. 추가 자료조사"

$ws.Cells.Item(14, 5).Value = "1. (2-2) 수익적 요소 추가`n2. (4-1) 기능정의 추가`n3. 추가 자료 조사`n - 비슷한 어플(마이클) 화면 구성 조사"
$ws.Cells.Item(14, 5).HorizontalAlignment = -4131
$ws.Cells.Item(14, 5).VerticalAlignment = -4108
$ws.Cells.Item(14, 5).WrapText = $true

# Row 15: add Goal (D15)
$ws.Cells.Item(15, 4).Value = "1. 추가 자료조사`n2. 요구사항 추가 조사`n3. 경쟁 어플 기능 비교 보완"
$ws.Cells.Item(15, 4).HorizontalAlignment = -4131
$ws.Cells.Item(15, 4).VerticalAlignment = -4108
$ws.Cells.Item(15, 4).WrapText = $true

$ws.Range("D15").Select() | Out-Null

Write-Host "done"
